$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("OR-1", "PV-001", 13),
    @("OR-2", "PV-001", 17),
    @("OR-2", "PV-001", 11),
    @("OR-2", "PV-002", 13),
    @("OR-3", "PV-001", 1),
    @("OR-3", "PV-003", 6),
    @("OR-3", "PV-003", 11),
    @("OR-4", "PV-002", 11),
    @("OR-5", "PV-002", 11),
    @("OR-6", "PV-003", 9),
    @("OR-6", "PV-003", 18),
    @("OR-6", "PV-001", 12),
    @("OR-7", "PV-001", 13),
    @("OR-7", "PV-001", 3),
    @("OR-8", "PV-001", 16),
    @("OR-9", "PV-002", 15),
    @("OR-9", "PV-002", 15),
    @("OR-9", "PV-001", 20),
    @("OR-10", "PV-001", 15),
    @("OR-10", "PV-003", 5),
    @("OR-11", "PV-003", 9),
    @("OR-11", "PV-003", 9),
    @("OR-12", "PV-003", 17),
    @("OR-12", "PV-002", 16),
    @("OR-12", "PV-003", 16),
    @("OR-13", "PV-002", 14),
    @("OR-13", "PV-002", 14),
    @("OR-13", "PV-001", 3),
    @("OR-14", "PV-001", 4),
    @("OR-15", "PV-003", 1),
    @("OR-15", "PV-003", 4),
    @("OR-16", "PV-002", 19),
    @("OR-16", "PV-001", 11),
    @("OR-17", "PV-002", 8),
    @("OR-17", "PV-003", 19),
    @("OR-17", "PV-001", 17),
    @("OR-18", "PV-001", 12),
    @("OR-18", "PV-003", 8),
    @("OR-19", "PV-003", 15),
    @("OR-19", "PV-001", 2),
    @("OR-20", "PV-002", 11),
    @("OR-21", "PV-001", 13),
    @("OR-21", "PV-002", 12),
    @("OR-21", "PV-001", 19),
    @("OR-22", "PV-003", 18),
    @("OR-22", "PV-001", 4),
    @("OR-22", "PV-001", 12),
    @("OR-23", "PV-001", 5),
    @("OR-24", "PV-003", 15),
    @("OR-25", "PV-003", 3),
    @("OR-25", "PV-001", 20),
    @("OR-26", "PV-001", 16),
    @("OR-27", "PV-003", 15),
    @("OR-27", "PV-003", 8),
    @("OR-27", "PV-002", 15),
    @("OR-28", "PV-001", 14),
    @("OR-28", "PV-002", 1),
    @("OR-29", "PV-002", 13),
    @("OR-29", "PV-001", 8),
    @("OR-30", "PV-001", 15),
    @("OR-30", "PV-001", 16),
    @("OR-31", "PV-001", 9),
    @("OR-32", "PV-003", 12)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $rowIndex++
}

$ws.Range("A2:C64").Select()
